$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRAMECALCULATOR")

# Bring the sheet into view scrolled down a bit (row 5 at the top) like the
# author left it - harmless if the host doesn't persist window scroll state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5

# Core edit: the cross-section size input dropped from 415 to 315. Every other
# changed cell in the workbook (H3, the summary string, the I8:I55 lookup
# totals, and the "BED Holes Location" Q column) is a formula that depends on
# this cell (directly or via FRAMECALCULATOR!H2/H3), so they recompute on
# their own once this value changes and the workbook recalculates.
$ws.Range("H2").Value = 315
